$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the MuSCs sending-cluster block (original rows 8-10), shifting remaining rows up.
# Excel automatically drops the now-unused "MuSCs" shared-string entry and renumbers the rest.
$ws.Range("A8:T10").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Update all remaining data rows (2-10) with recalculated TPM-derived figures.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl4"
$ws.Range("C2").Value = "Ccr1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2400906666666667
$ws.Range("H2").Value = 0.720272
$ws.Range("I2").Value = 0.001631540293869566
$ws.Range("J2").Value = 0.001631540293869566
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01219
$ws.Range("N2").Value = 0.03657
$ws.Range("O2").Value = 0.0002880357555630755
$ws.Range("P2").Value = 0.0002880357555630755
$ws.Range("Q2").Value = 0.002926705226666667
$ws.Range("R2").Value = 0.02634034704
$ws.Range("S2").Value = 0.0000004699419412763228
$ws.Range("T2").Value = 0.0000004699419412763227

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl4"
$ws.Range("C3").Value = "Ccr1"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2400906666666667
$ws.Range("H3").Value = 0.720272
$ws.Range("I3").Value = 0.001631540293869566
$ws.Range("J3").Value = 0.001631540293869566
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 32.87103466666667
$ws.Range("N3").Value = 98.613104
$ws.Range("O3").Value = 0.7767049471988007
$ws.Range("P3").Value = 0.7767049471988008
$ws.Range("Q3").Value = 7.892028627143111
$ws.Range("R3").Value = 71.028257644288
$ws.Range("S3").Value = 0.001267225417802677
$ws.Range("T3").Value = 0.001267225417802677

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl4"
$ws.Range("C4").Value = "Ccr1"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2400906666666667
$ws.Range("H4").Value = 0.720272
$ws.Range("I4").Value = 0.001631540293869566
$ws.Range("J4").Value = 0.001631540293869566
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.43791
$ws.Range("N4").Value = 28.31373
$ws.Range("O4").Value = 0.2230070170456362
$ws.Range("P4").Value = 0.2230070170456362
$ws.Range("Q4").Value = 2.26595410384
$ws.Range("R4").Value = 20.39358693456
$ws.Range("S4").Value = 0.0003638449341256127
$ws.Range("T4").Value = 0.0003638449341256127

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Ccl4"
$ws.Range("C5").Value = "Ccr1"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 91.82408133333333
$ws.Range("H5").Value = 275.472244
$ws.Range("I5").Value = 0.623992138981758
$ws.Range("J5").Value = 0.623992138981758
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01219
$ws.Range("N5").Value = 0.03657
$ws.Range("O5").Value = 0.0002880357555630755
$ws.Range("P5").Value = 0.0002880357555630755
$ws.Range("Q5").Value = 1.119335551453333
$ws.Range("R5").Value = 10.07401996308
$ws.Range("S5").Value = 0.0001797320472170303
$ws.Range("T5").Value = 0.0001797320472170303

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Ccl4"
$ws.Range("C6").Value = "Ccr1"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 91.82408133333333
$ws.Range("H6").Value = 275.472244
$ws.Range("I6").Value = 0.623992138981758
$ws.Range("J6").Value = 0.623992138981758
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 32.87103466666667
$ws.Range("N6").Value = 98.613104
$ws.Range("O6").Value = 0.7767049471988007
$ws.Range("P6").Value = 0.7767049471988008
$ws.Range("Q6").Value = 3018.352560742819
$ws.Range("R6").Value = 27165.17304668538
$ws.Range("S6").Value = 0.484657781360293
$ws.Range("T6").Value = 0.4846577813602931

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Ccl4"
$ws.Range("C7").Value = "Ccr1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 91.82408133333333
$ws.Range("H7").Value = 275.472244
$ws.Range("I7").Value = 0.623992138981758
$ws.Range("J7").Value = 0.623992138981758
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.43791
$ws.Range("N7").Value = 28.31373
$ws.Range("O7").Value = 0.2230070170456362
$ws.Range("P7").Value = 0.2230070170456362
$ws.Range("Q7").Value = 866.62741545668
$ws.Range("R7").Value = 7799.64673911012
$ws.Range("S7").Value = 0.1391546255742479
$ws.Range("T7").Value = 0.1391546255742479

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Ccl4"
$ws.Range("C8").Value = "Ccr1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 55.09165833333333
$ws.Range("H8").Value = 165.274975
$ws.Range("I8").Value = 0.3743763207243725
$ws.Range("J8").Value = 0.3743763207243724
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01219
$ws.Range("N8").Value = 0.03657
$ws.Range("O8").Value = 0.0002880357555630755
$ws.Range("P8").Value = 0.0002880357555630755
$ws.Range("Q8").Value = 0.6715673150833333
$ws.Range("R8").Value = 6.04410583575
$ws.Range("S8").Value = 0.0001078337664047689
$ws.Range("T8").Value = 0.0001078337664047689

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Ccl4"
$ws.Range("C9").Value = "Ccr1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 55.09165833333333
$ws.Range("H9").Value = 165.274975
$ws.Range("I9").Value = 0.3743763207243725
$ws.Range("J9").Value = 0.3743763207243724
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 32.87103466666667
$ws.Range("N9").Value = 98.613104
$ws.Range("O9").Value = 0.7767049471988007
$ws.Range("P9").Value = 0.7767049471988008
$ws.Range("Q9").Value = 1810.919810919156
$ws.Range("R9").Value = 16298.2782982724
$ws.Range("S9").Value = 0.290779940420705
$ws.Range("T9").Value = 0.290779940420705

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Ccl4"
$ws.Range("C10").Value = "Ccr1"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 55.09165833333333
$ws.Range("H10").Value = 165.274975
$ws.Range("I10").Value = 0.3743763207243725
$ws.Range("J10").Value = 0.3743763207243724
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.43791
$ws.Range("N10").Value = 28.31373
$ws.Range("O10").Value = 0.2230070170456362
$ws.Range("P10").Value = 0.2230070170456362
$ws.Range("Q10").Value = 519.9501131007501
$ws.Range("R10").Value = 4679.551017906751
$ws.Range("S10").Value = 0.08348854653726272
$ws.Range("T10").Value = 0.0834885465372627
